# Append 6 new data rows (159-164) to the "dataset" sheet, mirroring the
# existing tensorflow/ranking rows already present at the bottom of the
# table (rows 153-158). This matches the commit's "changes in file upload"
# which appended more rows for the same repository/commit date.

$excel = New-Object -ComObject Excel.Application
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the last 6 existing rows (values + formatting) straight down into the
# 6 new rows so the new rows inherit the same column formatting (bold id
# column, bordered cells, etc.) as the rest of the table.
$source = $ws.Range("A153:O158")
$target = $ws.Range("A159:O164")
$source.Copy($target)

# The copy duplicates row 158's "id" value (157) into every new row, so fix
# up column A with the correct sequential ids (158, 159, 160, 161, 162, 163).
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item(159 + $i, 1).Value = 158 + $i
}
